$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-set the date column format so new date cells adopt it directly
# instead of first picking up a default short-date format.
$ws.Range("B2:B33").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("A2").Value = "SI-1"
$ws.Range("B2").Value = [DateTime]"2025-01-01"
$ws.Range("C2").Value = "GOVERNMENT"
$ws.Range("D2").Value = "TEST-ID-14"
$ws.Range("E2").Value = "CUST-02"
$ws.Range("F2").Value = "CLOSED"

$ws.Range("A3").Value = "SI-2"
$ws.Range("B3").Value = [DateTime]"2025-01-02"
$ws.Range("C3").Value = "CORPORATE"
$ws.Range("D3").Value = "TEST-ID-14"
$ws.Range("E3").Value = "CUST-01"
$ws.Range("F3").Value = "CLOSED"

$ws.Range("A4").Value = "SI-3"
$ws.Range("B4").Value = [DateTime]"2025-01-03"
$ws.Range("C4").Value = "PATIENT"
$ws.Range("D4").Value = "TEST-ID-14"
$ws.Range("E4").Value = "CUST-03"
$ws.Range("F4").Value = "CLOSED"

$ws.Range("A5").Value = "SI-4"
$ws.Range("B5").Value = [DateTime]"2025-01-04"
$ws.Range("C5").Value = "PATIENT"
$ws.Range("D5").Value = "TEST-ID-15"
$ws.Range("E5").Value = "CUST-01"
$ws.Range("F5").Value = "CLOSED"

$ws.Range("A6").Value = "SI-5"
$ws.Range("B6").Value = [DateTime]"2025-01-05"
$ws.Range("C6").Value = "GOVERNMENT"
$ws.Range("D6").Value = "TEST-ID-15"
$ws.Range("E6").Value = "CUST-02"
$ws.Range("F6").Value = "CLOSED"

$ws.Range("A7").Value = "SI-6"
$ws.Range("B7").Value = [DateTime]"2025-01-06"
$ws.Range("C7").Value = "WHOLESALE"
$ws.Range("D7").Value = "TEST-ID-15"
$ws.Range("E7").Value = "CUST-03"
$ws.Range("F7").Value = "CLOSED"

$ws.Range("A8").Value = "SI-7"
$ws.Range("B8").Value = [DateTime]"2025-01-07"
$ws.Range("C8").Value = "PATIENT"
$ws.Range("D8").Value = "TEST-ID-14"
$ws.Range("E8").Value = "CUST-03"
$ws.Range("F8").Value = "CLOSED"

$ws.Range("A9").Value = "SI-8"
$ws.Range("B9").Value = [DateTime]"2025-01-08"
$ws.Range("C9").Value = "GOVERNMENT"
$ws.Range("D9").Value = "TEST-ID-15"
$ws.Range("E9").Value = "CUST-02"
$ws.Range("F9").Value = "CLOSED"

$ws.Range("A10").Value = "SI-9"
$ws.Range("B10").Value = [DateTime]"2025-01-09"
$ws.Range("C10").Value = "PATIENT"
$ws.Range("D10").Value = "TEST-ID-15"
$ws.Range("E10").Value = "CUST-01"
$ws.Range("F10").Value = "CLOSED"

$ws.Range("A11").Value = "SI-10"
$ws.Range("B11").Value = [DateTime]"2025-01-10"
$ws.Range("C11").Value = "PATIENT"
$ws.Range("D11").Value = "TEST-ID-15"
$ws.Range("E11").Value = "CUST-01"
$ws.Range("F11").Value = "CLOSED"

$ws.Range("A12").Value = "SI-11"
$ws.Range("B12").Value = [DateTime]"2025-01-11"
$ws.Range("C12").Value = "WHOLESALE"
$ws.Range("D12").Value = "TEST-ID-14"
$ws.Range("E12").Value = "CUST-03"
$ws.Range("F12").Value = "CLOSED"

$ws.Range("A13").Value = "SI-12"
$ws.Range("B13").Value = [DateTime]"2025-01-12"
$ws.Range("C13").Value = "PATIENT"
$ws.Range("D13").Value = "TEST-ID-14"
$ws.Range("E13").Value = "CUST-01"
$ws.Range("F13").Value = "CLOSED"

$ws.Range("A14").Value = "SI-13"
$ws.Range("B14").Value = [DateTime]"2025-01-13"
$ws.Range("C14").Value = "GOVERNMENT"
$ws.Range("D14").Value = "TEST-ID-14"
$ws.Range("E14").Value = "CUST-03"
$ws.Range("F14").Value = "CLOSED"

$ws.Range("A15").Value = "SI-14"
$ws.Range("B15").Value = [DateTime]"2025-01-14"
$ws.Range("C15").Value = "CORPORATE"
$ws.Range("D15").Value = "TEST-ID-15"
$ws.Range("E15").Value = "CUST-02"
$ws.Range("F15").Value = "CLOSED"

$ws.Range("A16").Value = "SI-15"
$ws.Range("B16").Value = [DateTime]"2025-01-15"
$ws.Range("C16").Value = "CORPORATE"
$ws.Range("D16").Value = "TEST-ID-14"
$ws.Range("E16").Value = "CUST-03"
$ws.Range("F16").Value = "CLOSED"

$ws.Range("A17").Value = "SI-16"
$ws.Range("B17").Value = [DateTime]"2025-01-16"
$ws.Range("C17").Value = "PATIENT"
$ws.Range("D17").Value = "TEST-ID-14"
$ws.Range("E17").Value = "CUST-02"
$ws.Range("F17").Value = "CLOSED"

$ws.Range("A18").Value = "SI-17"
$ws.Range("B18").Value = [DateTime]"2025-01-17"
$ws.Range("C18").Value = "PATIENT"
$ws.Range("D18").Value = "TEST-ID-15"
$ws.Range("E18").Value = "CUST-03"
$ws.Range("F18").Value = "CLOSED"

$ws.Range("A19").Value = "SI-18"
$ws.Range("B19").Value = [DateTime]"2025-01-18"
$ws.Range("C19").Value = "GOVERNMENT"
$ws.Range("D19").Value = "TEST-ID-14"
$ws.Range("E19").Value = "CUST-02"
$ws.Range("F19").Value = "CLOSED"

$ws.Range("A20").Value = "SI-19"
$ws.Range("B20").Value = [DateTime]"2025-01-19"
$ws.Range("C20").Value = "GOVERNMENT"
$ws.Range("D20").Value = "TEST-ID-15"
$ws.Range("E20").Value = "CUST-01"
$ws.Range("F20").Value = "CLOSED"

$ws.Range("A21").Value = "SI-20"
$ws.Range("B21").Value = [DateTime]"2025-01-20"
$ws.Range("C21").Value = "GOVERNMENT"
$ws.Range("D21").Value = "TEST-ID-14"
$ws.Range("E21").Value = "CUST-03"
$ws.Range("F21").Value = "CLOSED"

$ws.Range("A22").Value = "SI-21"
$ws.Range("B22").Value = [DateTime]"2025-01-21"
$ws.Range("C22").Value = "CORPORATE"
$ws.Range("D22").Value = "TEST-ID-15"
$ws.Range("E22").Value = "CUST-01"
$ws.Range("F22").Value = "CLOSED"

$ws.Range("A23").Value = "SI-22"
$ws.Range("B23").Value = [DateTime]"2025-01-22"
$ws.Range("C23").Value = "PATIENT"
$ws.Range("D23").Value = "TEST-ID-14"
$ws.Range("E23").Value = "CUST-02"
$ws.Range("F23").Value = "CLOSED"

$ws.Range("A24").Value = "SI-23"
$ws.Range("B24").Value = [DateTime]"2025-01-23"
$ws.Range("C24").Value = "GOVERNMENT"
$ws.Range("D24").Value = "TEST-ID-14"
$ws.Range("E24").Value = "CUST-02"
$ws.Range("F24").Value = "CLOSED"

$ws.Range("A25").Value = "SI-24"
$ws.Range("B25").Value = [DateTime]"2025-01-24"
$ws.Range("C25").Value = "PATIENT"
$ws.Range("D25").Value = "TEST-ID-14"
$ws.Range("E25").Value = "CUST-03"
$ws.Range("F25").Value = "CLOSED"

$ws.Range("A26").Value = "SI-25"
$ws.Range("B26").Value = [DateTime]"2025-01-25"
$ws.Range("C26").Value = "RETAIL"
$ws.Range("D26").Value = "TEST-ID-15"
$ws.Range("E26").Value = "CUST-01"
$ws.Range("F26").Value = "CLOSED"

$ws.Range("A27").Value = "SI-26"
$ws.Range("B27").Value = [DateTime]"2025-01-26"
$ws.Range("C27").Value = "CORPORATE"
$ws.Range("D27").Value = "TEST-ID-15"
$ws.Range("E27").Value = "CUST-01"
$ws.Range("F27").Value = "CLOSED"

$ws.Range("A28").Value = "SI-27"
$ws.Range("B28").Value = [DateTime]"2025-01-27"
$ws.Range("C28").Value = "GOVERNMENT"
$ws.Range("D28").Value = "TEST-ID-15"
$ws.Range("E28").Value = "CUST-01"
$ws.Range("F28").Value = "CLOSED"

$ws.Range("A29").Value = "SI-28"
$ws.Range("B29").Value = [DateTime]"2025-01-28"
$ws.Range("C29").Value = "RETAIL"
$ws.Range("D29").Value = "TEST-ID-14"
$ws.Range("E29").Value = "CUST-01"
$ws.Range("F29").Value = "CLOSED"

$ws.Range("A30").Value = "SI-29"
$ws.Range("B30").Value = [DateTime]"2025-01-29"
$ws.Range("C30").Value = "CORPORATE"
$ws.Range("D30").Value = "TEST-ID-14"
$ws.Range("E30").Value = "CUST-02"
$ws.Range("F30").Value = "CLOSED"

$ws.Range("A31").Value = "SI-30"
$ws.Range("B31").Value = [DateTime]"2025-01-30"
$ws.Range("C31").Value = "WHOLESALE"
$ws.Range("D31").Value = "TEST-ID-14"
$ws.Range("E31").Value = "CUST-02"
$ws.Range("F31").Value = "CLOSED"

$ws.Range("A32").Value = "SI-31"
$ws.Range("B32").Value = [DateTime]"2025-01-31"
$ws.Range("C32").Value = "CORPORATE"
$ws.Range("D32").Value = "TEST-ID-15"
$ws.Range("E32").Value = "CUST-03"
$ws.Range("F32").Value = "CLOSED"

$ws.Range("A33").Value = "SI-32"
$ws.Range("B33").Value = [DateTime]"2025-02-01"
$ws.Range("C33").Value = "GOVERNMENT"
$ws.Range("D33").Value = "TEST-ID-14"
$ws.Range("E33").Value = "CUST-03"
$ws.Range("F33").Value = "CLOSED"

$ws.Range("A2:F33").Select()
